# Applies scheduled-runner value updates across multiple worksheets in the
# Golem Profits workbook. Each worksheet corresponds to a crafting job
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) and the rows below hold market
# pricing data in columns H-N that get refreshed by the scheduled runner.

$wb = $excel.ActiveWorkbook

function Set-RowValues($SheetName, $Row, $Values) {
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in $Values.Keys) {
        $cell = $col + $Row
        if ($null -eq $Values[$col]) {
            $ws.Range($cell).ClearContents()
        } else {
            $ws.Range($cell).Value = $Values[$col]
        }
    }
}

# ALC row 9
Set-RowValues "ALC" 9 @{
    H = 793.6
    I = 42.75
    J = 3797
    K = 42.75
    L = 3797
    M = 126.25
    N = -4135
}

# ALC row 34
Set-RowValues "ALC" 34 @{
    H = 49333.332
    I = 8000
    K = 8000
    M = -7797
}

# ALC row 36
Set-RowValues "ALC" 36 @{
    H = 49333.332
    I = 8000
    K = 8000
    M = -7285
}

# ALC row 40
Set-RowValues "ALC" 40 @{
    H = 2259.4
    I = 2099.6667
    K = 2099.6667
    M = -1924.6667
}

# ALC row 54
Set-RowValues "ALC" 54 @{
    H = 15000
    I = 0
    J = 15000
    K = 0
    L = $null
    M = 15000
    N = -15972
}

# ALC row 116
Set-RowValues "ALC" 116 @{
    H = 4398.2856
    I = 4598.5
    K = 4598.5
    M = -1156.5
}

# ALC row 132
Set-RowValues "ALC" 132 @{
    H = 67866.10000000001
    I = 67866.10000000001
    K = 203598.3
    M = -201068.3
}

# ALC row 137
Set-RowValues "ALC" 137 @{
    H = 1166.8182
    I = 1022.55554
    K = 3067.66662
    M = -517.66662
}

# ARM row 2
Set-RowValues "ARM" 2 @{
    H = 766.7778
    I = 766.7778
    K = 766.7778
    M = -653.7778
}

# ARM row 116
Set-RowValues "ARM" 116 @{
    H = 766.7778
    I = 766.7778
    K = 766.7778
    M = 1527.2222
}

# ARM row 122
Set-RowValues "ARM" 122 @{
    H = 1290
    I = 1290
    K = 3870
    M = -1420
}

# ARM row 131
Set-RowValues "ARM" 131 @{
    H = 49000
    J = 49000
    L = 49000
    N = -59080
}

# ARM row 132
Set-RowValues "ARM" 132 @{
    H = 2048.5
    I = 2048.5
    J = 0
    K = 6145.5
    L = 0
    M = -3615.5
    N = $null
}

# BSM row 3
Set-RowValues "BSM" 3 @{
    H = 766.7778
    I = 766.7778
    K = 766.7778
    M = -652.7778
}

# BSM row 134
Set-RowValues "BSM" 134 @{
    H = 0
    I = 0
    J = 0
    K = 0
    L = 0
    M = $null
    N = 0
}

# CRP row 16
Set-RowValues "CRP" 16 @{
    H = 676.5
    I = 553
    J = 800
    K = 553
    L = 800
    M = -266
    N = -1374
}

# CRP row 113
Set-RowValues "CRP" 113 @{
    H = 676.5
    I = 553
    J = 800
    K = 553
    L = 800
    M = 1617
    N = -5140
}

# CRP row 132
Set-RowValues "CRP" 132 @{
    H = 2252.25
    I = 2336.3333
    J = 2000
    K = 7008.999899999999
    L = 6000
    M = -4478.999899999999
    N = -11060
}

# CUL row 80
Set-RowValues "CUL" 80 @{
    H = 5898
    I = 5898
    J = 0
    K = 17694
    L = 0
    M = $null
    N = -16758
}

# CUL row 83
Set-RowValues "CUL" 83 @{
    H = 5898
    I = 5898
    J = 0
    K = 53082
    L = 0
    M = $null
    N = -48402
}

# CUL row 86
Set-RowValues "CUL" 86 @{
    H = 380.2
    J = 412.25
    L = 1236.75
    N = -3608.75
}

# CUL row 89
Set-RowValues "CUL" 89 @{
    H = 380.2
    J = 412.25
    L = 3710.25
    N = -15566.25
}

# GSM row 70
Set-RowValues "GSM" 70 @{
    H = 100003140
    I = 3733.3333
    K = 3733.3333
    M = -3463.3333
}

# GSM row 73
Set-RowValues "GSM" 73 @{
    H = 100003140
    I = 3733.3333
    K = 3733.3333
    M = -2797.3333
}

# GSM row 126
Set-RowValues "GSM" 126 @{
    H = 1999.75
    I = 1933
    J = 2200
    K = 5799
    L = 6600
    M = -3329
    N = -11540
}

# LTW row 22
Set-RowValues "LTW" 22 @{
    H = 3296.923
    J = 3481.2856
    L = 3481.2856
    N = -4071.2856
}

# LTW row 27
Set-RowValues "LTW" 27 @{
    H = 3296.923
    J = 3481.2856
    L = 3481.2856
    N = -3695.2856
}

# LTW row 39
Set-RowValues "LTW" 39 @{
    H = 23205
    J = 34912.5
    L = 34912.5
    N = -35832.5
}

# LTW row 46
Set-RowValues "LTW" 46 @{
    H = 1650
    J = 2500
    L = 2500
    N = -2876
}

# LTW row 93
Set-RowValues "LTW" 93 @{
    H = 55557040
    I = 83335064
    J = 993
    K = 83335064
    L = 993
    M = -83333816
    N = -3489
}

# LTW row 122
Set-RowValues "LTW" 122 @{
    H = 3545.4
    I = 2663
    J = 4133.6665
    K = 7989
    L = 12400.9995
    M = -5539
    N = -17300.9995
}

# LTW row 132
Set-RowValues "LTW" 132 @{
    H = 1156.2858
    I = 1156.2858
    K = 3468.8574
    M = -938.8574000000003
}

# WVR row 126
Set-RowValues "WVR" 126 @{
    H = 4972.706
    I = 4205.5454
    K = 12616.6362
    M = -10146.6362
}

# WVR row 132
Set-RowValues "WVR" 132 @{
    H = 2914.2856
    I = 1800
    K = 5400
    M = -2870
}
